# Remove form_id from remaining forms
# (mirrors the commit: drop the "form_id" column on the settings sheet and
# shift the surrounding cell comments + simplify the conditional-format ranges
# on the survey sheet so they no longer reference the now-removed column.)

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")
$survey = $wb.Worksheets.Item("survey")

# --- settings sheet: drop the form_id column (column B) ------------------
# Deleting the whole column shifts version/style/namespaces left by one and
# updates the shared-string table / formula cell automatically.
$settings.Range("B1").EntireColumn.Delete()

# --- settings sheet: comments don't auto-shift with the column, so re-point
#     them by hand to describe what is now in each column.
$settings.Range("B1").Comment.Text("The unique version code that identifies the current state of the form. A common convention is to use a format like yyyymmddrr. For example, 2017021501 is the 1st revision from Feb 15th, 2017.`n`nBy default, this template uses a formula to create a date-based version that will update automatically.")
$settings.Range("C1").Comment.Text("Set to ‘pages’ to indicate that groups with the ``field-list`` appearance represent separate form pages (and all other questions will be shown on their own page). ")
$settings.Range("D1").Comment.Delete()
$settings.Range("E1").Comment.Delete()
$settings.Range("D1").AddComment("Custom namespaces supported in the form.  ``cht`` must be included here to use the custom ``instance::cht`` columns on the survey sheet.") | Out-Null

# Leave the selection on the settings sheet where the deleted column used to
# be, then restore "survey" as the active tab.
$settings.Range("B1").Select()
$survey.Activate()

# --- survey sheet: consolidate the conditional-formatting ranges ----------
# These sqref lists used to be split around the (now nonexistent) row 26
# carve-out; collapse each back into one contiguous range.
$fcWhole = $survey.Range("A2:G25").FormatConditions.Item(1)
$fcWhole.ModifyAppliesToRange($survey.Range("A2:G9999"))

$fcColC = $survey.Range("C2:C25").FormatConditions.Item(6)
$fcColC.ModifyAppliesToRange($survey.Range("C2:C9999"))
